# Update the spicule type list on Sheet1.
# Row 2 used to hold a single "SHOW PICTURE" placeholder (plus empty
# cells in B2:F2). It is replaced with three genuine data rows that
# list the spicule types recorded for this measurement sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "SHOW PICTURE" placeholder text in A2.
$ws.Range("A2").Value = "small Oxea"

# The old row 2 also had empty placeholder cells in B2:F2 - clear them
# out since the new layout only uses column A for these rows.
$ws.Range("B2:F2").ClearContents()

# Add the two additional spicule type rows.
$ws.Range("A3").Value = "large Oxea"
$ws.Range("A4").Value = "spined oxyaster euaster"

# Match the author's final selection/view state.
[void]$ws.Range("C6").Select()
